$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 and C2 updated
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 1

# Rename the three "plural" item labels to their singular form.
# (written in this order so the shared-strings table gets the new
# entries appended in the same sequence as the target workbook)
$ws.Range("A11").Value = "SteelOxygenBlownConverter"
$ws.Range("A16").Value = "SteelOpenHearthFurnace"
$ws.Range("A12").Value = "SteelElectricFurnace"

# All rows 11:37 get the same Minimum/Maximum ratio constraints now.
$ws.Range("B11:B37").Value = 0.5
$ws.Range("C11:C37").Value = 0.15

# Move the active selection to A11.
$ws.Range("A11").Select()
